$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the first three data rows (plus header) survive - remove the
# duplicated trailing rows 5-7.
$ws.Range("A5:A7").EntireRow.Delete()

# Row 2 (still size "S" / "Computer") now shows a different price.
$ws.Cells.Item(2, 3).Value = "Rs. 389"

# Row 4 used to be a duplicate of the "S" / Computer / Rs. 899 row; it now
# holds a brand new product record. Set column B first so the new shared
# string for the long product name is interned before the "M" string.
$ws.Cells.Item(4, 2).Value = "Muscle Torque Zip Sweatshirts"
$ws.Cells.Item(4, 1).Value = "M"
$ws.Cells.Item(4, 3).Value = "Rs. 1390"

# Column B needs to be widened so the long product name fits.
$ws.Columns.Item(2).AutoFit()

# Record the new active cell/selection.
[void]$ws.Range("O14").Select()

# Make sure the sheet is set up for portrait printing.
$ws.PageSetup.Orientation = 1

Write-Host "Applied TestData edits"
